$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New inventory rows (Ref GNA80496, dated 2025-09-25) appended below the
# existing data (rows 393-400), mirroring the upload adding a new batch
# of feather-cost entries broken out by size.
$newRows = @(
    @(393, 15, 92,  1380, "40X40"),
    @(394, 15, 215, 3225, "45X70"),
    @(395, 50, 114, 5700, "50X50"),
    @(396, 45, 132, 5940, "55X55"),
    @(397, 60, 161, 9660, "60X60"),
    @(398, 45, 198, 8910, "65X65"),
    @(399, 10, 241, 2410, "70X70"),
    @(400, 30, 93,  2790, "35X55")
)

foreach ($row in $newRows) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = "GNA80496"
    $ws.Cells.Item($r, 2).Value = 45925
    $ws.Cells.Item($r, 3).Value = $row[1]
    $ws.Cells.Item($r, 4).Value = $row[2]
    $ws.Cells.Item($r, 5).Value = $row[3]
    $ws.Cells.Item($r, 6).Value = $row[4]
}

# Scroll/selection state matches the saved view after entering the new data.
$ws.Range("F400").Select() | Out-Null
